$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-06 Saturday", "2024-04-07 Sunday"),
    @("335÷3=", "650÷7="),
    @("569÷6=", "623÷7="),
    @("673÷3=", "143÷3="),
    @("853÷5=", "290÷8="),
    @("277÷9=", "674÷6="),
    @("953÷7=", "320÷6="),
    @("326÷2=", "199÷4="),
    @("563÷6=", "388÷4="),
    @("122÷3=", "301÷8="),
    @("441÷6=", "872÷3="),
    @("626÷8=", "436÷7="),
    @("918÷4=", "860÷3="),
    @("305÷8=", "702÷2="),
    @("971÷5=", "121÷8="),
    @("452÷2=", "584÷8="),
    @("775÷5=", "662÷4="),
    @("116÷3=", "604÷2="),
    @("308÷8=", "434÷7="),
    @("314÷7=", "630÷5="),
    @("896÷5=", "113÷7="),
    @("173÷6=", "288÷6="),
    @("119÷3=", "602÷8="),
    @("379÷4=", "632÷3="),
    @("460÷4=", "397÷2="),
    @("421÷5=", "620÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
